{"js": "// Update the date line and the division problems in the table, matching\n// the new day's generated worksheet content. Replacements are applied by\n// paragraph index (document order) so that duplicate expressions (e.g.\n// \"701\u00f78=\" appears twice in the original) are each replaced with the\n// correct, distinct new value.\nconst replacements = [\n  [0, \"2025-07-27 Sunday\"],\n  [1, \"396\u00f76=\"],\n  [2, \"595\u00f79=\"],\n  [3, \"883\u00f74=\"],\n  [4, \"764\u00f76=\"],\n  [5, \"245\u00f77=\"],\n  [21, \"896\u00f78=\"],\n  [22, \"813\u00f73=\"],\n  [23, \"298\u00f76=\"],\n  [24, \"712\u00f76=\"],\n  [25, \"246\u00f74=\"],\n  [41, \"285\u00f78=\"],\n  [42, \"790\u00f72=\"],\n  [43, \"622\u00f73=\"],\n  [44, \"969\u00f73=\"],\n  [45, \"730\u00f79=\"],\n  [61, \"809\u00f75=\"],\n  [62, \"414\u00f77=\"],\n  [63, \"972\u00f79=\"],\n  [64, \"366\u00f74=\"],\n  [65, \"740\u00f79=\"],\n  [81, \"240\u00f73=\"],\n  [82, \"764\u00f74=\"],\n  [83, \"408\u00f75=\"],\n  [84, \"913\u00f73=\"],\n  [85, \"759\u00f78=\"],\n];\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst maxIndex = replacements.reduce((m, [i]) => Math.max(m, i), 0);\nif (paragraphs.items.length <= maxIndex) {\n  throw new Error(\n    \"Expected at least \" + (maxIndex + 1) + \" paragraphs, found \" + paragraphs.items.length\n  );\n}\n\nfor (const [index, newText] of replacements) {\n  // insertText with \"Replace\" keeps the existing run formatting (font,\n  // size, etc.) of the paragraph while swapping in the new text.\n  paragraphs.items[index].insertText(newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Update the date line and the division problems in the table, matching\n# the new day's generated worksheet content.\n$d = $word.ActiveDocument\n\n# --- Date paragraph (first paragraph in the document body, above table) ---\n$dateRange = $d.Paragraphs.Item(1).Range\n$dateRange.MoveEnd(1, -1)   # wdCharacter = 1; exclude the paragraph mark\n$dateRange.Text = \"2025-07-27 Sunday\"\n\n# --- Division problems table ---\n# The table has 20 rows x 5 columns; only every 4th row (1, 5, 9, 13, 17)\n# actually holds content, the rest are blank spacer rows.\n$tbl = $d.Tables.Item(1)\n\n$newValues = @{\n    1  = @(\"396\u00f76=\", \"595\u00f79=\", \"883\u00f74=\", \"764\u00f76=\", \"245\u00f77=\")\n    5  = @(\"896\u00f78=\", \"813\u00f73=\", \"298\u00f76=\", \"712\u00f76=\", \"246\u00f74=\")\n    9  = @(\"285\u00f78=\", \"790\u00f72=\", \"622\u00f73=\", \"969\u00f73=\", \"730\u00f79=\")\n    13 = @(\"809\u00f75=\", \"414\u00f77=\", \"972\u00f79=\", \"366\u00f74=\", \"740\u00f79=\")\n    17 = @(\"240\u00f73=\", \"764\u00f74=\", \"408\u00f75=\", \"913\u00f73=\", \"759\u00f78=\")\n}\n\nforeach ($rowIndex in @(1, 5, 9, 13, 17)) {\n    $values = $newValues[$rowIndex]\n    for ($col = 1; $col -le 5; $col++) {\n        $cell = $tbl.Cell($rowIndex, $col)\n        $cell.Range.Text = $values[$col - 1]\n    }\n}\n"}
